$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the previous logbook entry block (rows 7:8, merged B7:J8) down to
# rows 9:10 so the new entry inherits identical borders/formatting.
$src = $ws.Range("A7:J8")
$dst = $ws.Range("A9:J10")
$src.Copy($dst)

# Fill in the new entry's date and comment.
$ws.Range("A9").Value = 43167
$ws.Range("B9").Value = "Added SENSE_R global variable to sense resistor circuit. Changed BJTs to MOSFET equivalents. Also added text, descriptions for different blocks of circuit. "

# Match the author's final selection.
$ws.Range("I14").Select() | Out-Null
